$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsOverview.Range("G2").Value = "2016-08-15 12:59:18"

$wsZhCn.Range("H2").Value = "2016-08-15 12:59:14"
$wsZhCn.Range("K2").Value = "2016-08-15 12:59:32"

$wsDeDe.Range("H2").Value = "2016-08-15 12:59:18"
$wsDeDe.Range("K2").Value = "2016-08-15 12:59:38"
